$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Vulnerabilities = Server vulnerabilities + OS vulnerabilities`r") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Server vulnerabilities = " + [char]0x2026
        break
    }
}
